$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9747682213783264
$ws.Range("B1").Value = 1.022755980491638
$ws.Range("C1").Value = 6.494430541992188
$ws.Range("D1").Value = 1.968039035797119
$ws.Range("E1").Value = 1.106608867645264
